$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "39.699.33", "6.00").
# Excel's Value setter auto-coerces such literals to floating-point numbers,
# which silently drops trailing zeros / loses exact decimal text. Forcing the
# cell to Text format for the duration of the write keeps the literal string,
# then resetting the style back to Normal avoids leaving a stray number format.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '39.699.33'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.172.57'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.09%  '
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('E10').Value = '  -1.43%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('E12').Value = '  -1.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.493.91'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.811'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.49'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.37%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.174.45'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '39.648.01'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0910'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.73%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.58'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.01%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('E29').Value = '  +1.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.83'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.67'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.94%  '
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.53'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.39%  '
$ws.Range('E34').Value = '  -2.34%  '
$ws.Range('E35').Value = '  -2.52%  '
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.87'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.59%  '
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('E40').Value = '  +13.08%  '
$ws.Range('E41').Value = '  -1.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '102.58'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.74'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.513.32'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E45').Value = '  +2.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.89'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0920'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.79%  '
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('E50').Value = '  +33.37%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '49.48'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.31%  '
